$d = $word.ActiveDocument

# --- Paragraph 8 (the "relative abundance" paragraph) -----------------
$old8 = "This analysis centers around determining where healthcare and service workers respectively are relatively abundant. I’ve done this by constructing a relative abundance index that takes into account both the density of the type of worker in a give census block group, and the proportion of people in that block group that work in a given sector. In this way the relative abundance measure incorporates not only how many workers reside in a given block group, but also highlights where those workers are more likely to live."

$new8 = "This analysis centers around determining where healthcare and service workers respectively are relatively abundant. To do this the data utilizes ArcGIS’s hot spot analysis tool, which uses the Getis-Ord Gi* statistic to determine areas of high and low values using the context of the surrounding areas. By highlighting areas with more essential workers, we can see where the analysis should further focus."

$rng8 = $d.Content
$found8 = $rng8.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)

# --- Paragraph 9 (the "bus usage" / ratio paragraph) -------------------
$old9 = "Third, this project takes in bus usage data from the time of the bus changes, 3-25-20, until 4-21-20, and compares that data to the change in bus frequency. For this analysis I’ve constructed a ratio that compares that change in use to the change in buses—values above 1 means use has decreased by less than the bus availability has decreased, or that use has actually increased. Values below 1 represent where buses were reduced more than use declined, and these values should receive scrutiny."

$new9 = "Third, this project takes in bus usage data from the time of the bus changes, 3-25-20, until 4-21-20, and compares that data to the change in bus frequency. For this analysis I’ve constructed a ratio that compares that change in use to the change in buses—values above 1 means use has decreased by less than the bus availability has decreased, or that use has actually increased and should receive scrutiny. Many stops with values above 1.0 saw less than one rider per day on average, so the analysis highlights stops that see more than 50 riders daily on average to bring to the Port Authority’s attention. The story map explaining the findings can be found at https://arcg.is/1q1XC4."

$rng9 = $d.Content
$found9 = $rng9.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)

Write-Output "found8=$found8 found9=$found9"
